$d = $word.ActiveDocument

# --- Locate the existing OVERVIEW specifications table (the Property/Value
# header-only table) and the three legacy Jinja loop paragraphs that
# immediately follow it ({% for spec %} / {{ spec.property }}: {{ spec.value }} / {% endfor %}) ---

$loopStartPara = $null
$loopEndPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($loopStartPara -eq $null -and $t -like "*{% for spec in overview_specifications_table %}*") {
        $loopStartPara = $p
    }
    if ($t -like "*{% endfor %}*") {
        $loopEndPara = $p
    }
}

if ($loopStartPara -eq $null -or $loopEndPara -eq $null) {
    throw "Could not locate the overview_specifications_table loop paragraphs"
}

$insertionPoint = $loopStartPara.Range.Start

# --- Build the new two-column specification table (header row + 8 data
# rows bound to overview_specifications_table[N].value) as raw WordprocessingML
# and insert it immediately before the loop paragraphs ---

$tableXml = '<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:type="auto" w:w="0"/><w:tblLook w:firstColumn="1" w:firstRow="1" w:lastColumn="0" w:lastRow="0" w:noHBand="0" w:noVBand="1" w:val="04A0"/></w:tblPr><w:tblGrid><w:gridCol w:w="5400"/><w:gridCol w:w="5400"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Property</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Value</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Product Name</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[0].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Reactive Species</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[1].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Size</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[2].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[3].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Sensitivity</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[4].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Detection Range</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[5].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Storage Instructions</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[6].value }}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>Uniprot ID</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:type="dxa" w:w="5400"/></w:tcPr><w:p><w:r><w:t>{{ overview_specifications_table[7].value }}</w:t></w:r></w:p></w:tc></w:tr></w:tbl>'

$insertRange = $d.Range($insertionPoint, $insertionPoint)
$insertRange.InsertXML($tableXml)

# --- Remove the now-redundant Jinja loop paragraphs. Re-resolve the range
# fresh (the document shifted after InsertXML) by searching again. ---

$loopStartPara = $null
$loopEndPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($loopStartPara -eq $null -and $t -like "*{% for spec in overview_specifications_table %}*") {
        $loopStartPara = $p
    }
    if ($t -like "*{% endfor %}*") {
        $loopEndPara = $p
    }
}

$deleteRange = $d.Range($loopStartPara.Range.Start, $loopEndPara.Range.End)
$deleteRange.Delete()

Write-Host "Overview specifications table generated."
